$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.461.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.08%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.572.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.01%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -0.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'286.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.75%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3652"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.85%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -3.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.3339"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.81%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.128"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.62%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07443"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.39%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.06%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'20.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.75%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.988"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.85%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.54%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.575.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.22%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001110"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.18%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'88.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.06%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06740"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.30%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.396"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.63%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'16.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.69%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'12.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.52%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'22.451.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.09%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.386"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.40%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.617"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.28%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'152.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.17%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.82%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.021"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.52%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'124.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.83%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.751.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.26%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.036"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.16%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'6.168"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.61%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.994"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.01%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'9.766"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.38%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.08273"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.96%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.02431"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.96%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.2263"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.66%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.06477"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.416"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.48%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.300"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.03%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'11.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.54%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.6272"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.68%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'14.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.68%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.6067"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.79%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.67%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.048"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.02%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'124.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.69%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.222"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.81%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.07223"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.42%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'76.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.18%  "
$ws.Range("E51").Style = "Normal"
